$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 810
$ws.Range("F7").Value = 380
$ws.Range("F8").Value = 4631
$ws.Range("F9").Value = 4631
$ws.Range("F11").Value = 115
$ws.Range("F12").Value = 150
$ws.Range("F15").Value = 104
$ws.Range("F16").Value = 7287
$ws.Range("F18").Value = 124
$ws.Range("F21").Value = 504
$ws.Range("F22").Value = 1319
$ws.Range("F28").Value = 6132
$ws.Range("F31").Value = 111
$ws.Range("F34").Value = 6328
$ws.Range("F37").Value = 94
$ws.Range("F46").Value = 408
$ws.Range("F47").Value = 2116
$ws.Range("F49").Value = 1064

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G5").Value = 118
$ws.Range("F6").Value = 115

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G6").Value = 118
$ws.Range("F7").Value = 380
$ws.Range("F8").Value = 4631
$ws.Range("F9").Value = 4631
$ws.Range("F11").Value = 115
$ws.Range("F12").Value = 150
$ws.Range("F15").Value = 104
$ws.Range("F16").Value = 7287
$ws.Range("F18").Value = 124
$ws.Range("F19").Value = 504
$ws.Range("F20").Value = 1319
$ws.Range("F21").Value = 115
$ws.Range("F28").Value = 6132
$ws.Range("F32").Value = 111
$ws.Range("F35").Value = 6328
$ws.Range("F38").Value = 94
$ws.Range("F46").Value = 408
$ws.Range("F48").Value = 2116
